$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new value for column F ("dSF")
$updates = @{
    2  = -4
    4  = 1
    5  = -2
    6  = 0
    7  = 1
    8  = -1
    9  = -4
    10 = -6
    11 = -1
    12 = -6
    13 = -5
    14 = -6
    15 = -3
    16 = -1
    18 = -3
    19 = -8
    20 = 4
    22 = -4
    23 = 2
    24 = -4
    25 = -3
    26 = 1
    27 = -3
    29 = -4
    30 = 3
    31 = -5
    33 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
